$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage values that must land as plain shared-string
# text (no style, no quote-prefix) even though they look numeric - Excel's
# normal Value assignment auto-converts "03250001"-style strings to numbers
# and stripping the leading zeros, so we build the text via a formula result
# and paste-special the *value only* onto the real cell, which keeps the
# destination at the default style (matches cells typed as text originally).
$scratch = "Z100"

function Set-TextValue([string]$addr, [string]$text) {
    $escaped = $text.Replace("""", """""")
    $ws.Range($scratch).Formula = "=""" + $escaped + """"
    $ws.Range($scratch).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range($scratch).Clear()
}

# ------------------------------------------------------------------
# Header row (row 1) - 14 columns (A:N), inserting two new headers
# after "Date de creation" (Date de cloture, Durée de résolution)
# and one new header before "description" (Cloturé par).
# ------------------------------------------------------------------
$ws.Range("A1").Value = "NumRef"
$ws.Range("B1").Value = "Date de creation"
$ws.Range("C1").Value = "Date de cloture"
$ws.Range("D1").Value = "Durée de résolution"
$ws.Range("E1").Value = "Type d'incident"
$ws.Range("F1").Value = "Cause d'incident"
$ws.Range("G1").Value = "Equipement"
$ws.Range("H1").Value = "Site"
$ws.Range("I1").Value = "Shift"
$ws.Range("J1").Value = "Utilisateur"
$ws.Range("K1").Value = "Cloturé par"
$ws.Range("L1").Value = "description"
$ws.Range("M1").Value = "Édité par"
$ws.Range("N1").Value = "Status"

# ------------------------------------------------------------------
# Row 2 - existing incident record, now carrying the two extra date /
# duration columns plus the new "Cloturé par" column, and its Status
# moved from "EN ATTENTE" to "EN MAINTENANCE".
# ------------------------------------------------------------------
Set-TextValue "A2" "03250001"
$ws.Range("B2").Value = 45720.35496090278

# C2 is a new closing-date cell - clone the date number format from B2
# (style index 1) via copy/paste-special so no new style gets minted.
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = 45727.67790376158

$ws.Range("D2").Value = "175 Heure(s)"
$ws.Range("E2").Value = "COUPURE ELECTRIQUE"
$ws.Range("F2").Value = "L'AXE LOURD ÉTAIT BLOQUE A CAUSE DE LA CAN, CAR DES JOUEURS SORTAIENT POUR LES ENTRAINEMENTS."
$ws.Range("G2").Value = "AIRE DE PESEES"
$ws.Range("H2").Value = "SALLE D'ATTENTE NIVEAU 1"
$ws.Range("I2").Value = "Shift(14-22)"
$ws.Range("J2").Value = "Admin User"
$ws.Range("K2").Value = "N/A"
$ws.Range("L2").Value = "dfgdf"
$ws.Range("N2").Value = "EN MAINTENANCE"

# ------------------------------------------------------------------
# Row 3 - brand new incident record.
# ------------------------------------------------------------------
Set-TextValue "A3" "02254432"

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = 45693.61901905092

$ws.Range("D3").Value = "-482990 Heure(s)"
$ws.Range("E3").Value = "COUPURE ELECTRIQUE"
$ws.Range("F3").Value = "COUPURE ELECTRIQUE"
$ws.Range("G3").Value = "GROUPE ELECTROGENE"
$ws.Range("H3").Value = "P02"
$ws.Range("J3").Value = "FOTSO TSOBGNY FRANCK JOEL"
$ws.Range("K3").Value = "N/A"
$ws.Range("L3").Value = "On a eu coupure d'électricité sur notre pont"
$ws.Range("N3").Value = "CLOTURE"

# ------------------------------------------------------------------
# Column widths - shift the old layout two columns to the right from
# column C on, matching the widened table.
# ------------------------------------------------------------------
$w15 = 15 - 5/6
$w20 = 20 - 5/6
$w50 = 50 - 5/6

$ws.Columns("A").ColumnWidth = $w15
$ws.Columns("B:C").ColumnWidth = $w20
$ws.Columns("D:F").ColumnWidth = $w50
$ws.Columns("G").ColumnWidth = $w15
$ws.Columns("H:K").ColumnWidth = $w20
$ws.Columns("L").ColumnWidth = $w50
$ws.Columns("M:N").ColumnWidth = $w20
